$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely - this shifts all rows below (3..22) up by one
# position, so the data that was in row 3 becomes row 2, etc., and the
# last row (22) is removed, shrinking the used range to A1:M21.
$ws.Rows.Item(2).Delete()
